$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 202.25  # ALC H5: 202.44444 -> 202.25
$ws.Cells.Item(5, 9).Value = 228.28572  # ALC I5: 202.44444 -> 228.28572
$ws.Cells.Item(5, 10).Value = 20  # ALC J5: 0 -> 20
$ws.Cells.Item(5, 11).Value = 228.28572  # ALC K5: 202.44444 -> 228.28572
$ws.Cells.Item(5, 12).Value = 20  # ALC L5: 0 -> 20
$ws.Cells.Item(5, 14).Value = -250  # ALC N5: None -> -250

$ws.Cells.Item(12, 8).Value = 731  # ALC H12: 799.2 -> 731
$ws.Cells.Item(12, 9).Value = 797.2  # ALC I12: 899 -> 797.2
$ws.Cells.Item(12, 11).Value = 797.2  # ALC K12: 899 -> 797.2
$ws.Cells.Item(12, 13).Value = -627.2  # ALC M12: -729 -> -627.2

$ws.Cells.Item(20, 8).Value = 3999  # ALC H20: 4000 -> 3999
$ws.Cells.Item(20, 9).Value = 3999  # ALC I20: 4000 -> 3999
$ws.Cells.Item(20, 11).Value = 3999  # ALC K20: 4000 -> 3999
$ws.Cells.Item(20, 13).Value = -3769  # ALC M20: -3770 -> -3769

$ws.Cells.Item(35, 8).Value = 3999  # ALC H35: 4000 -> 3999
$ws.Cells.Item(35, 9).Value = 3999  # ALC I35: 4000 -> 3999
$ws.Cells.Item(35, 11).Value = 3999  # ALC K35: 4000 -> 3999
$ws.Cells.Item(35, 13).Value = -3620  # ALC M35: -3621 -> -3620

$ws.Cells.Item(46, 8).Value = 3000  # ALC H46: 4000 -> 3000
$ws.Cells.Item(46, 10).Value = 3000  # ALC J46: 4000 -> 3000
$ws.Cells.Item(46, 12).Value = 9000  # ALC L46: 12000 -> 9000
$ws.Cells.Item(46, 14).Value = -9238  # ALC N46: -12238 -> -9238

$ws.Cells.Item(53, 8).Value = 193.66667  # ALC H53: 253.125 -> 193.66667
$ws.Cells.Item(53, 9).Value = 241  # ALC I53: 331.83334 -> 241
$ws.Cells.Item(53, 10).Value = 99  # ALC J53: 17 -> 99
$ws.Cells.Item(53, 11).Value = 241  # ALC K53: 331.83334 -> 241
$ws.Cells.Item(53, 12).Value = 99  # ALC L53: 17 -> 99
$ws.Cells.Item(53, 13).Value = 396  # ALC M53: 305.16666 -> 396
$ws.Cells.Item(53, 14).Value = -1373  # ALC N53: -1291 -> -1373

$ws.Cells.Item(60, 8).Value = 3000  # ALC H60: 4000 -> 3000
$ws.Cells.Item(60, 10).Value = 3000  # ALC J60: 4000 -> 3000
$ws.Cells.Item(60, 12).Value = 9000  # ALC L60: 12000 -> 9000
$ws.Cells.Item(60, 14).Value = -9968  # ALC N60: -12968 -> -9968

$ws.Cells.Item(100, 8).Value = 2836.2727  # ALC H100: 3155.7778 -> 2836.2727
$ws.Cells.Item(100, 9).Value = 2800.125  # ALC I100: 3071.8572 -> 2800.125
$ws.Cells.Item(100, 10).Value = 2932.6667  # ALC J100: 3449.5 -> 2932.6667
$ws.Cells.Item(100, 11).Value = 2800.125  # ALC K100: 3071.8572 -> 2800.125
$ws.Cells.Item(100, 12).Value = 2932.6667  # ALC L100: 3449.5 -> 2932.6667
$ws.Cells.Item(100, 13).Value = -2259.125  # ALC M100: -2530.8572 -> -2259.125
$ws.Cells.Item(100, 14).Value = -4014.6667  # ALC N100: -4531.5 -> -4014.6667

$ws.Cells.Item(112, 8).Value = 3185.7778  # ALC H112: 3423.3 -> 3185.7778
$ws.Cells.Item(112, 10).Value = 3171.75  # ALC J112: 3454.625 -> 3171.75
$ws.Cells.Item(112, 12).Value = 9515.25  # ALC L112: 10363.875 -> 9515.25
$ws.Cells.Item(112, 14).Value = -11731.25  # ALC N112: -12579.875 -> -11731.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 29644  # ARM H76: 29312.666 -> 29644
$ws.Cells.Item(76, 10).Value = 29644  # ARM J76: 29312.666 -> 29644
$ws.Cells.Item(76, 12).Value = 29644  # ARM L76: 29312.666 -> 29644
$ws.Cells.Item(76, 14).Value = -30320  # ARM N76: -29988.666 -> -30320

$ws.Cells.Item(79, 8).Value = 29644  # ARM H79: 29312.666 -> 29644
$ws.Cells.Item(79, 10).Value = 29644  # ARM J79: 29312.666 -> 29644
$ws.Cells.Item(79, 12).Value = 29644  # ARM L79: 29312.666 -> 29644
$ws.Cells.Item(79, 14).Value = -31984  # ARM N79: -31652.666 -> -31984

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2970.1875  # BSM H86: 3251.9285 -> 2970.1875
$ws.Cells.Item(86, 9).Value = 3101.6  # BSM I86: 3425.2307 -> 3101.6
$ws.Cells.Item(86, 11).Value = 3101.6  # BSM K86: 3425.2307 -> 3101.6
$ws.Cells.Item(86, 13).Value = -1978.6  # BSM M86: -2302.2307 -> -1978.6

$ws.Cells.Item(89, 8).Value = 2970.1875  # BSM H89: 3251.9285 -> 2970.1875
$ws.Cells.Item(89, 9).Value = 3101.6  # BSM I89: 3425.2307 -> 3101.6
$ws.Cells.Item(89, 11).Value = 15508  # BSM K89: 17126.1535 -> 15508
$ws.Cells.Item(89, 13).Value = -9892  # BSM M89: -11510.1535 -> -9892

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2998.5  # CRP H99: 2999 -> 2998.5
$ws.Cells.Item(99, 9).Value = 2998.5  # CRP I99: 2999 -> 2998.5
$ws.Cells.Item(99, 11).Value = 2998.5  # CRP K99: 2999 -> 2998.5
$ws.Cells.Item(99, 13).Value = -1500.5  # CRP M99: -1501 -> -1500.5

$ws.Cells.Item(126, 8).Value = 2998.5  # CRP H126: 2999 -> 2998.5
$ws.Cells.Item(126, 9).Value = 2998.5  # CRP I126: 2999 -> 2998.5
$ws.Cells.Item(126, 11).Value = 8995.5  # CRP K126: 8997 -> 8995.5
$ws.Cells.Item(126, 13).Value = -6525.5  # CRP M126: -6527 -> -6525.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 994.7895  # CUL H12: 996.8421 -> 994.7895
$ws.Cells.Item(12, 9).Value = 155.25  # CUL I12: 177.28572 -> 155.25
$ws.Cells.Item(12, 10).Value = 1605.3636  # CUL J12: 1474.9166 -> 1605.3636
$ws.Cells.Item(12, 11).Value = 465.75  # CUL K12: 531.85716 -> 465.75
$ws.Cells.Item(12, 12).Value = 4816.0908  # CUL L12: 4424.7498 -> 4816.0908
$ws.Cells.Item(12, 13).Value = -292.75  # CUL M12: -358.85716 -> -292.75
$ws.Cells.Item(12, 14).Value = -5162.0908  # CUL N12: -4770.7498 -> -5162.0908

$ws.Cells.Item(15, 8).Value = 424.7143  # CUL H15: 382.25 -> 424.7143
$ws.Cells.Item(15, 9).Value = 259.6  # CUL I15: 230.5 -> 259.6
$ws.Cells.Item(15, 11).Value = 778.8000000000001  # CUL K15: 691.5 -> 778.8000000000001
$ws.Cells.Item(15, 13).Value = -638.8000000000001  # CUL M15: -551.5 -> -638.8000000000001

$ws.Cells.Item(60, 8).Value = 3085.2222  # CUL H60: 3295.2856 -> 3085.2222
$ws.Cells.Item(60, 9).Value = 2395.6667  # CUL I60: 2487 -> 2395.6667
$ws.Cells.Item(60, 11).Value = 7187.000100000001  # CUL K60: 7461 -> 7187.000100000001
$ws.Cells.Item(60, 13).Value = -6936.000100000001  # CUL M60: -7210 -> -6936.000100000001

$ws.Cells.Item(62, 8).Value = 8548.5  # CUL H62: 7864.3335 -> 8548.5
$ws.Cells.Item(62, 10).Value = 6497  # CUL J62: 6496.5 -> 6497
$ws.Cells.Item(62, 12).Value = 19491  # CUL L62: 19489.5 -> 19491
$ws.Cells.Item(62, 14).Value = -20863  # CUL N62: -20861.5 -> -20863

$ws.Cells.Item(65, 8).Value = 8548.5  # CUL H65: 7864.3335 -> 8548.5
$ws.Cells.Item(65, 10).Value = 6497  # CUL J65: 6496.5 -> 6497
$ws.Cells.Item(65, 12).Value = 58473  # CUL L65: 58468.5 -> 58473
$ws.Cells.Item(65, 14).Value = -65337  # CUL N65: -65332.5 -> -65337

$ws.Cells.Item(75, 8).Value = 1751.25  # CUL H75: 1584.1666 -> 1751.25
$ws.Cells.Item(75, 9).Value = 1503  # CUL I75: 1502 -> 1503
$ws.Cells.Item(75, 10).Value = 1999.5  # CUL J75: 1666.3334 -> 1999.5
$ws.Cells.Item(75, 11).Value = 4509  # CUL K75: 4506 -> 4509
$ws.Cells.Item(75, 12).Value = 5998.5  # CUL L75: 4999.0002 -> 5998.5
$ws.Cells.Item(75, 13).Value = -3511  # CUL M75: -3508 -> -3511
$ws.Cells.Item(75, 14).Value = -7994.5  # CUL N75: -6995.0002 -> -7994.5

$ws.Cells.Item(78, 8).Value = 1751.25  # CUL H78: 1584.1666 -> 1751.25
$ws.Cells.Item(78, 9).Value = 1503  # CUL I78: 1502 -> 1503
$ws.Cells.Item(78, 10).Value = 1999.5  # CUL J78: 1666.3334 -> 1999.5
$ws.Cells.Item(78, 11).Value = 13527  # CUL K78: 13518 -> 13527
$ws.Cells.Item(78, 12).Value = 17995.5  # CUL L78: 14997.0006 -> 17995.5
$ws.Cells.Item(78, 13).Value = -8535  # CUL M78: -8526 -> -8535
$ws.Cells.Item(78, 14).Value = -27979.5  # CUL N78: -24981.0006 -> -27979.5

$ws.Cells.Item(107, 8).Value = 705.36365  # CUL H107: 700.8182 -> 705.36365
$ws.Cells.Item(107, 9).Value = 676.6  # CUL I107: 663.6667 -> 676.6
$ws.Cells.Item(107, 10).Value = 729.3333  # CUL J107: 745.4 -> 729.3333
$ws.Cells.Item(107, 11).Value = 2029.8  # CUL K107: 1991.0001 -> 2029.8
$ws.Cells.Item(107, 12).Value = 2187.9999  # CUL L107: 2236.2 -> 2187.9999
$ws.Cells.Item(107, 13).Value = -109.8000000000002  # CUL M107: -71.00009999999997 -> -109.8000000000002
$ws.Cells.Item(107, 14).Value = -6027.9999  # CUL N107: -6076.2 -> -6027.9999

$ws.Cells.Item(132, 8).Value = 1356.2858  # CUL H132: 1265.8334 -> 1356.2858
$ws.Cells.Item(132, 10).Value = 1899.5  # CUL J132: 1900 -> 1899.5
$ws.Cells.Item(132, 12).Value = 17095.5  # CUL L132: 17100 -> 17095.5
$ws.Cells.Item(132, 14).Value = -22155.5  # CUL N132: -22160 -> -22155.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 9000  # GSM H22: 7000 -> 9000
$ws.Cells.Item(22, 10).Value = 9000  # GSM J22: 7000 -> 9000
$ws.Cells.Item(22, 12).Value = 9000  # GSM L22: 7000 -> 9000
$ws.Cells.Item(22, 14).Value = -10058  # GSM N22: -8058 -> -10058

$ws.Cells.Item(75, 8).Value = 50000  # GSM H75: 0 -> 50000
$ws.Cells.Item(75, 9).Value = 50000  # GSM I75: 0 -> 50000
$ws.Cells.Item(75, 11).Value = 50000  # GSM K75: 0 -> 50000
$ws.Cells.Item(75, 13).Value = -49126  # GSM M75: None -> -49126

$ws.Cells.Item(78, 8).Value = 50000  # GSM H78: 0 -> 50000
$ws.Cells.Item(78, 9).Value = 50000  # GSM I78: 0 -> 50000
$ws.Cells.Item(78, 11).Value = 150000  # GSM K78: 0 -> 150000
$ws.Cells.Item(78, 13).Value = -145632  # GSM M78: None -> -145632

$ws.Cells.Item(102, 8).Value = 2026.7646  # GSM H102: 1888.5 -> 2026.7646
$ws.Cells.Item(102, 9).Value = 1965.9375  # GSM I102: 1823.1177 -> 1965.9375
$ws.Cells.Item(102, 11).Value = 1965.9375  # GSM K102: 1823.1177 -> 1965.9375
$ws.Cells.Item(102, 13).Value = -343.9375  # GSM M102: -201.1177 -> -343.9375

$ws.Cells.Item(122, 8).Value = 3983.5715  # GSM H122: 4223 -> 3983.5715
$ws.Cells.Item(122, 9).Value = 3358.875  # GSM I122: 3629.5 -> 3358.875
$ws.Cells.Item(122, 11).Value = 10076.625  # GSM K122: 10888.5 -> 10076.625
$ws.Cells.Item(122, 13).Value = -7626.625  # GSM M122: -8438.5 -> -7626.625

$ws.Cells.Item(123, 8).Value = 34210.5  # GSM H123: 50000 -> 34210.5
$ws.Cells.Item(123, 10).Value = 34210.5  # GSM J123: 50000 -> 34210.5
$ws.Cells.Item(123, 12).Value = 34210.5  # GSM L123: 50000 -> 34210.5
$ws.Cells.Item(123, 14).Value = -39110.5  # GSM N123: -54900 -> -39110.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2900  # LTW H22: 2350 -> 2900
$ws.Cells.Item(22, 9).Value = 0  # LTW I22: 1800 -> 0
$ws.Cells.Item(22, 11).Value = 0  # LTW K22: 1800 -> 0
$ws.Cells.Item(22, 13).ClearContents()  # LTW M22: remove (was -1505)

$ws.Cells.Item(27, 8).Value = 2900  # LTW H27: 2350 -> 2900
$ws.Cells.Item(27, 9).Value = 0  # LTW I27: 1800 -> 0
$ws.Cells.Item(27, 11).Value = 0  # LTW K27: 1800 -> 0
$ws.Cells.Item(27, 13).ClearContents()  # LTW M27: remove (was -1693)

$ws.Cells.Item(30, 8).Value = 1057.5  # LTW H30: 1027.6666 -> 1057.5
$ws.Cells.Item(30, 9).Value = 1169  # LTW I30: 1163.2 -> 1169
$ws.Cells.Item(30, 10).Value = 500  # LTW J30: 350 -> 500
$ws.Cells.Item(30, 11).Value = 1169  # LTW K30: 1163.2 -> 1169
$ws.Cells.Item(30, 12).Value = 500  # LTW L30: 350 -> 500
$ws.Cells.Item(30, 13).Value = -1061  # LTW M30: -1055.2 -> -1061
$ws.Cells.Item(30, 14).Value = -716  # LTW N30: -566 -> -716

$ws.Cells.Item(40, 8).Value = 8138  # LTW H40: 8585.571 -> 8138
$ws.Cells.Item(40, 9).Value = 8110.125  # LTW I40: 8720 -> 8110.125
$ws.Cells.Item(40, 11).Value = 8110.125  # LTW K40: 8720 -> 8110.125
$ws.Cells.Item(40, 13).Value = -7974.125  # LTW M40: -8584 -> -7974.125

$ws.Cells.Item(55, 8).Value = 851.38464  # LTW H55: 791.2857 -> 851.38464
$ws.Cells.Item(55, 9).Value = 709.8570999999999  # LTW I55: 623.625 -> 709.8570999999999
$ws.Cells.Item(55, 10).Value = 1016.5  # LTW J55: 1014.8333 -> 1016.5
$ws.Cells.Item(55, 11).Value = 709.8570999999999  # LTW K55: 623.625 -> 709.8570999999999
$ws.Cells.Item(55, 12).Value = 1016.5  # LTW L55: 1014.8333 -> 1016.5
$ws.Cells.Item(55, 13).Value = -536.8570999999999  # LTW M55: -450.625 -> -536.8570999999999
$ws.Cells.Item(55, 14).Value = -1362.5  # LTW N55: -1360.8333 -> -1362.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1214  # WVR H100: 1248 -> 1214
$ws.Cells.Item(100, 9).Value = 321  # WVR I100: 330.66666 -> 321
$ws.Cells.Item(100, 10).Value = 3000  # WVR J100: 4000 -> 3000
$ws.Cells.Item(100, 11).Value = 642  # WVR K100: 661.33332 -> 642
$ws.Cells.Item(100, 12).Value = 6000  # WVR L100: 8000 -> 6000
$ws.Cells.Item(100, 13).Value = -101  # WVR M100: -120.33332 -> -101
$ws.Cells.Item(100, 14).Value = -7082  # WVR N100: -9082 -> -7082

$ws.Cells.Item(122, 8).Value = 2375.05  # WVR H122: 2422.4736 -> 2375.05
$ws.Cells.Item(122, 9).Value = 2636.1333  # WVR I122: 2719.1428 -> 2636.1333
$ws.Cells.Item(122, 11).Value = 7908.3999  # WVR K122: 8157.428400000001 -> 7908.3999
$ws.Cells.Item(122, 13).Value = -5458.3999  # WVR M122: -5707.428400000001 -> -5458.3999

$ws.Cells.Item(132, 8).Value = 3786.5264  # WVR H132: 4139.353 -> 3786.5264
$ws.Cells.Item(132, 9).Value = 3061.7058  # WVR I132: 3364.9333 -> 3061.7058
$ws.Cells.Item(132, 11).Value = 9185.117400000001  # WVR K132: 10094.7999 -> 9185.117400000001
$ws.Cells.Item(132, 13).Value = -6655.117400000001  # WVR M132: -7564.7999 -> -6655.117400000001

$ws.Cells.Item(135, 8).Value = 112000  # WVR H135: 92999.664 -> 112000
$ws.Cells.Item(135, 9).Value = 0  # WVR I135: 54999 -> 0
$ws.Cells.Item(135, 11).Value = 0  # WVR K135: 54999 -> 0
$ws.Cells.Item(135, 13).ClearContents()  # WVR M135: remove (was -49929)
